$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 46 - shifts existing rows 46:149 down to 47:150
$ws.Rows(46).Insert()

# Populate the newly inserted row 46 with the new price-record data
$ws.Range("A46").Value = 7
$ws.Range("B46").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C46").Value = "Ñuble"
$ws.Range("D46").Value = 44519
$ws.Range("E46").Value = 16
$ws.Range("F46").Value = 100112003
$ws.Range("G46").Value = "Ajo"
$ws.Range("H46").Value = "Chino"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 100
$ws.Range("K46").Value = 18000
$ws.Range("L46").Value = 19000
$ws.Range("M46").Value = 18500
$ws.Range("N46").Value = "`$/caja 10 kilos"
$ws.Range("O46").Value = "China"
$ws.Range("P46").Value = 1850
$ws.Range("Q46").Value = 10
$ws.Range("R46").Value = "Hortaliza"
